$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = "sv"
$ws.Range("J5").Value = "Statement-opinion"
$ws.Range("I11").Value = "sv"
$ws.Range("J11").Value = "Statement-opinion"
$ws.Range("I12").Value = "aa"
$ws.Range("J12").Value = "Agree/Accept"
$ws.Range("I15").Value = "sd"
$ws.Range("J15").Value = "Statement-non-opinion"
$ws.Range("I38").Value = "sd"
$ws.Range("J38").Value = "Statement-non-opinion"
$ws.Range("I41").Value = "sv"
$ws.Range("J41").Value = "Statement-opinion"
$ws.Range("I46").Value = "aa"
$ws.Range("J46").Value = "Agree/Accept"
$ws.Range("I53").Value = "%"
$ws.Range("J53").Value = "Uninterpretable"
$ws.Range("I58").Value = "sd"
$ws.Range("J58").Value = "Statement-non-opinion"
$ws.Range("I60").Value = "b"
$ws.Range("J60").Value = "Acknowledge (Backchannel)"
$ws.Range("I77").Value = "aa"
$ws.Range("J77").Value = "Agree/Accept"
$ws.Range("I83").Value = "ba"
$ws.Range("J83").Value = "Appreciation"
$ws.Range("I88").Value = "sd"
$ws.Range("J88").Value = "Statement-non-opinion"
$ws.Range("I89").Value = "sd"
$ws.Range("J89").Value = "Statement-non-opinion"
$ws.Range("I113").Value = "sv"
$ws.Range("J113").Value = "Statement-opinion"
$ws.Range("I115").Value = "sd"
$ws.Range("J115").Value = "Statement-non-opinion"
$ws.Range("I118").Value = "aa"
$ws.Range("J118").Value = "Agree/Accept"
$ws.Range("I121").Value = "sv"
$ws.Range("J121").Value = "Statement-opinion"
$ws.Range("I129").Value = "aa"
$ws.Range("J129").Value = "Agree/Accept"
$ws.Range("I131").Value = "sv"
$ws.Range("J131").Value = "Statement-opinion"
$ws.Range("I134").Value = "aa"
$ws.Range("J134").Value = "Agree/Accept"
$ws.Range("I135").Value = "b"
$ws.Range("J135").Value = "Acknowledge (Backchannel)"
$ws.Range("I148").Value = "sv"
$ws.Range("J148").Value = "Statement-opinion"
$ws.Range("I153").Value = "sd"
$ws.Range("J153").Value = "Statement-non-opinion"
$ws.Range("I154").Value = "aa"
$ws.Range("J154").Value = "Agree/Accept"
$ws.Range("I156").Value = "sd"
$ws.Range("J156").Value = "Statement-non-opinion"
$ws.Range("I171").Value = "ba"
$ws.Range("J171").Value = "Appreciation"
$ws.Range("I196").Value = "aa"
$ws.Range("J196").Value = "Agree/Accept"
$ws.Range("I202").Value = "sv"
$ws.Range("J202").Value = "Statement-opinion"
$ws.Range("I205").Value = "ba"
$ws.Range("J205").Value = "Appreciation"
$ws.Range("I214").Value = "sv"
$ws.Range("J214").Value = "Statement-opinion"
$ws.Range("I215").Value = "sd"
$ws.Range("J215").Value = "Statement-non-opinion"
$ws.Range("I220").Value = "sv"
$ws.Range("J220").Value = "Statement-opinion"
$ws.Range("I230").Value = "aa"
$ws.Range("J230").Value = "Agree/Accept"
$ws.Range("I246").Value = "sd"
$ws.Range("J246").Value = "Statement-non-opinion"
$ws.Range("I260").Value = "ba"
$ws.Range("J260").Value = "Appreciation"
$ws.Range("I278").Value = "sd"
$ws.Range("J278").Value = "Statement-non-opinion"
$ws.Range("I309").Value = "sv"
$ws.Range("J309").Value = "Statement-opinion"
$ws.Range("I310").Value = "sd"
$ws.Range("J310").Value = "Statement-non-opinion"
$ws.Range("I327").Value = "%"
$ws.Range("J327").Value = "Uninterpretable"
$ws.Range("I328").Value = "%"
$ws.Range("J328").Value = "Uninterpretable"
$ws.Range("I342").Value = "sd"
$ws.Range("J342").Value = "Statement-non-opinion"
$ws.Range("I344").Value = "sd"
$ws.Range("J344").Value = "Statement-non-opinion"
$ws.Range("I383").Value = "sd"
$ws.Range("J383").Value = "Statement-non-opinion"
$ws.Range("I394").Value = "ba"
$ws.Range("J394").Value = "Appreciation"
$ws.Range("I402").Value = "sv"
$ws.Range("J402").Value = "Statement-opinion"
$ws.Range("I417").Value = "sv"
$ws.Range("J417").Value = "Statement-opinion"
$ws.Range("I421").Value = "aa"
$ws.Range("J421").Value = "Agree/Accept"
$ws.Range("I426").Value = "aa"
$ws.Range("J426").Value = "Agree/Accept"
$ws.Range("I438").Value = "sv"
$ws.Range("J438").Value = "Statement-opinion"
$ws.Range("I439").Value = "sv"
$ws.Range("J439").Value = "Statement-opinion"
$ws.Range("I441").Value = "sv"
$ws.Range("J441").Value = "Statement-opinion"
$ws.Range("I443").Value = "sd"
$ws.Range("J443").Value = "Statement-non-opinion"
$ws.Range("I457").Value = "sv"
$ws.Range("J457").Value = "Statement-opinion"
$ws.Range("I458").Value = "sd"
$ws.Range("J458").Value = "Statement-non-opinion"
$ws.Range("I461").Value = "sv"
$ws.Range("J461").Value = "Statement-opinion"
$ws.Range("I483").Value = "sv"
$ws.Range("J483").Value = "Statement-opinion"
$ws.Range("I484").Value = "ba"
$ws.Range("J484").Value = "Appreciation"
